$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = 1.512
$ws.Range("C13").Value = -3.64
$ws.Range("C39").Value = -2.137
$ws.Range("C45").Value = -3.35
$ws.Range("C57").Value = -3.445
$ws.Range("C66").Value = -3.737
$ws.Range("C76").Value = -4.672
$ws.Range("C82").Value = -3.167
$ws.Range("C90").Value = -2.431
$ws.Range("C94").Value = 3.765
$ws.Range("C97").Value = -0.491
